$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PR")

# --- Header: PR Number ---
$ws.Range("C7").Value = "PR No.:  2020-02-00102"

# --- Line item 1 (row 11) ---
$ws.Range("A11").Value = "S217"
$ws.Range("B11").Value = "pack"
$ws.Range("C11").Value = "Conduct of CFLGA-RIMTF Table Top Assessement (Live-in)" + [char]10 + "."
$ws.Range("D11").Value = 18
$ws.Range("E11").Value = 2000
$ws.Range("F11").Value = 36000

# --- Line item 2 (row 12) ---
$ws.Range("A12").Value = "S218"
$ws.Range("B12").Value = "pack"
$ws.Range("C12").Value = "Conduct of CFLGA-RIMTF Table Top Assessement (Live-out)" + [char]10 + "."
$ws.Range("D12").Value = 18
$ws.Range("E12").Value = 1200
$ws.Range("F12").Value = 21600

# --- Purpose ---
$ws.Range("B37").Value = "CFLGA  ASSESSMENT 2020"

# --- Approved by: printed name & designation ---
$ws.Range("D43").Value = "NOEL R. BARTOLABAC, CESO V"
$ws.Range("D44").Value = "Assistant Regional Director"

# --- Active cell selection ---
$ws.Range("D43").Select()
